$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (automated data refresh).
# All of these source cells are plain text in the workbook (prices are
# formatted strings, not numbers), so force a text number format before
# assigning values that would otherwise be auto-converted to numbers.
$ws.Range("D2").Value = "54.177.14"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.270.92"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "499.06"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.84"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0951"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.70"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").Value = "2.670.97"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.66"
$ws.Range("E14").Value = "  +4.56%  "
$ws.Range("D15").Value = "54.152.37"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "2.288.66"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.15"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "303.29"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.31"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.10"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.150"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.28"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.54"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0684"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.91"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.73"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.957"
$ws.Range("E34").Value = "  +10.60%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.82"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.77"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.544"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "238.45"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.76"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.18"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("E51").Value = "  -0.27%  "
